# Apply "Well Builder" updates to sheet "3. Well Builder"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3. Well Builder")

# Mark rows 17-20 column A with an "X" (matches existing pattern used throughout the workbook)
$ws.Range("A17").Value = "X"
$ws.Range("A18").Value = "X"
$ws.Range("A19").Value = "X"
$ws.Range("A20").Value = "X"

# New task rows 22-27
$ws.Range("A22").Value = "X"
$ws.Range("B22").Value = 13
$ws.Range("C22").Value = "Generate tubing info database"

$ws.Range("B23").Value = 14
$ws.Range("C23").Value = "Liner pulls from tubdata nd csgdata for sizes"

$ws.Range("B24").Value = 15
$ws.Range("C24").Value = "Casing info on left side"

$ws.Range("B25").Value = 16
$ws.Range("C25").Value = "formations on right side"

$ws.Range("B26").Value = 17
$ws.Range("C26").Value = "company, well name, and api on top"

$ws.Range("B27").Value = 18
$ws.Range("C27").Value = "remove casing button"

# Update the selection to reflect where the user left off editing
$ws.Range("B28").Select()
